$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 380.06668  # H6: was 387.5
$ws.Cells.Item(6, 9).Value = 284.69232  # I6: was 331.81818
$ws.Cells.Item(6, 11).Value = 854.07696  # K6: was 995.45454
$ws.Cells.Item(6, 13).Value = -742.07696  # M6: was -883.45454
$ws.Cells.Item(19, 8).Value = 2220.6553  # H19: was 2079.0967
$ws.Cells.Item(19, 9).Value = 4834.364  # I19: was 4435.75
$ws.Cells.Item(19, 10).Value = 623.3889  # J19: was 590.6842
$ws.Cells.Item(19, 11).Value = 4834.364  # K19: was 4435.75
$ws.Cells.Item(19, 12).Value = 623.3889  # L19: was 590.6842
$ws.Cells.Item(19, 13).Value = -4659.364  # M19: was -4260.75
$ws.Cells.Item(19, 14).Value = -973.3889  # N19: was -940.6842
$ws.Cells.Item(100, 8).Value = 4205.212  # H100: was 4728.593
$ws.Cells.Item(100, 9).Value = 1731.0834  # I100: was 1780.3636
$ws.Cells.Item(100, 10).Value = 5619  # J100: was 6755.5
$ws.Cells.Item(100, 11).Value = 1731.0834  # K100: was 1780.3636
$ws.Cells.Item(100, 12).Value = 5619  # L100: was 6755.5
$ws.Cells.Item(100, 13).Value = -1190.0834  # M100: was -1239.3636
$ws.Cells.Item(100, 14).Value = -6701  # N100: was -7837.5
$ws.Cells.Item(116, 8).Value = 154254.42  # H116: was 98725.55
$ws.Cells.Item(116, 9).Value = 179133.5  # I116: was 108100.1
$ws.Cells.Item(116, 11).Value = 179133.5  # K116: was 108100.1
$ws.Cells.Item(116, 13).Value = -175691.5  # M116: was -104658.1

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1177.6428  # H45: was 980.5599999999999
$ws.Cells.Item(45, 9).Value = 1060.875  # I45: was 915.65
$ws.Cells.Item(45, 10).Value = 1333.3334  # J45: was 1240.2
$ws.Cells.Item(45, 11).Value = 1060.875  # K45: was 915.65
$ws.Cells.Item(45, 12).Value = 1333.3334  # L45: was 1240.2
$ws.Cells.Item(45, 13).Value = -683.875  # M45: was -538.65
$ws.Cells.Item(45, 14).Value = -2087.3334  # N45: was -1994.2
$ws.Cells.Item(63, 8).Value = 1750  # H63: was 1800
$ws.Cells.Item(63, 9).Value = 1750  # I63: was 1800
$ws.Cells.Item(63, 11).Value = 1750  # K63: was 1800
$ws.Cells.Item(63, 13).Value = -1064  # M63: was -1114
$ws.Cells.Item(66, 8).Value = 1750  # H66: was 1800
$ws.Cells.Item(66, 9).Value = 1750  # I66: was 1800
$ws.Cells.Item(66, 11).Value = 8750  # K66: was 9000
$ws.Cells.Item(66, 13).Value = -5318  # M66: was -5568
$ws.Cells.Item(110, 8).Value = 1694.75  # H110: was 1937.9474
$ws.Cells.Item(110, 9).Value = 1670.1904  # I110: was 1951.3125
$ws.Cells.Item(110, 11).Value = 1670.1904  # K110: was 1951.3125
$ws.Cells.Item(110, 13).Value = 374.8096  # M110: was 93.6875
$ws.Cells.Item(122, 8).Value = 2056.818  # H122: was 1532.6333
$ws.Cells.Item(122, 9).Value = 1999.5883  # I122: was 1281.2727
$ws.Cells.Item(122, 10).Value = 2251.4  # J122: was 2223.875
$ws.Cells.Item(122, 11).Value = 5998.7649  # K122: was 3843.8181
$ws.Cells.Item(122, 12).Value = 6754.200000000001  # L122: was 6671.625
$ws.Cells.Item(122, 13).Value = -3548.7649  # M122: was -1393.8181
$ws.Cells.Item(122, 14).Value = -11654.2  # N122: was -11571.625
$ws.Cells.Item(132, 8).Value = 19406.467  # H132: was 32699.656
$ws.Cells.Item(132, 9).Value = 33032.727  # I132: was 76745.5
$ws.Cells.Item(132, 10).Value = 2752.1482  # J132: was 3335.762
$ws.Cells.Item(132, 11).Value = 99098.181  # K132: was 230236.5
$ws.Cells.Item(132, 12).Value = 8256.444600000001  # L132: was 10007.286
$ws.Cells.Item(132, 13).Value = -96568.181  # M132: was -227706.5
$ws.Cells.Item(132, 14).Value = -13316.4446  # N132: was -15067.286

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2052.423  # H31: was 2103.984
$ws.Cells.Item(31, 9).Value = 1107.409  # I31: was 1058.1702
$ws.Cells.Item(31, 10).Value = 7250  # J31: was 5380.8667
$ws.Cells.Item(31, 11).Value = 1107.409  # K31: was 1058.1702
$ws.Cells.Item(31, 12).Value = 7250  # L31: was 5380.8667
$ws.Cells.Item(31, 13).Value = -812.4090000000001  # M31: was -763.1702
$ws.Cells.Item(31, 14).Value = -7840  # N31: was -5970.8667
$ws.Cells.Item(34, 8).Value = 2052.423  # H34: was 2103.984
$ws.Cells.Item(34, 9).Value = 1107.409  # I34: was 1058.1702
$ws.Cells.Item(34, 10).Value = 7250  # J34: was 5380.8667
$ws.Cells.Item(34, 11).Value = 1107.409  # K34: was 1058.1702
$ws.Cells.Item(34, 12).Value = 7250  # L34: was 5380.8667
$ws.Cells.Item(34, 13).Value = -905.4090000000001  # M34: was -856.1702
$ws.Cells.Item(34, 14).Value = -7654  # N34: was -5784.8667
$ws.Cells.Item(58, 8).Value = 3883.4524  # H58: was 4247.1055
$ws.Cells.Item(58, 9).Value = 5425.952  # I58: was 6262.778
$ws.Cells.Item(58, 10).Value = 2340.9524  # J58: was 2433
$ws.Cells.Item(58, 11).Value = 5425.952  # K58: was 6262.778
$ws.Cells.Item(58, 12).Value = 2340.9524  # L58: was 2433
$ws.Cells.Item(58, 13).Value = -5222.952  # M58: was -6059.778
$ws.Cells.Item(58, 14).Value = -2746.9524  # N58: was -2839
$ws.Cells.Item(122, 8).Value = 1078.1765  # H122: was 3243.2222
$ws.Cells.Item(122, 9).Value = 1083.6875  # I122: was 3243.2222
$ws.Cells.Item(122, 10).Value = 990  # J122: was 0
$ws.Cells.Item(122, 11).Value = 3251.0625  # K122: was 9729.6666
$ws.Cells.Item(122, 12).Value = 2970  # L122: was 0
$ws.Cells.Item(122, 13).Value = -801.0625  # M122: was -7279.6666
$ws.Cells.Item(122, 14).Value = -7870  # N122: was None
$ws.Cells.Item(132, 8).Value = 2190.182  # H132: was 2397.1853
$ws.Cells.Item(132, 9).Value = 1514.9474  # I132: was 1583.25
$ws.Cells.Item(132, 10).Value = 3106.5715  # J132: was 3581.0908
$ws.Cells.Item(132, 11).Value = 4544.8422  # K132: was 4749.75
$ws.Cells.Item(132, 12).Value = 9319.7145  # L132: was 10743.2724
$ws.Cells.Item(132, 13).Value = -2014.8422  # M132: was -2219.75
$ws.Cells.Item(132, 14).Value = -14379.7145  # N132: was -15803.2724
$ws.Cells.Item(136, 8).Value = 3883.4524  # H136: was 4247.1055
$ws.Cells.Item(136, 9).Value = 5425.952  # I136: was 6262.778
$ws.Cells.Item(136, 10).Value = 2340.9524  # J136: was 2433
$ws.Cells.Item(136, 11).Value = 16277.856  # K136: was 18788.334
$ws.Cells.Item(136, 12).Value = 7022.8572  # L136: was 7299
$ws.Cells.Item(136, 13).Value = -13727.856  # M136: was -16238.334
$ws.Cells.Item(136, 14).Value = -12122.8572  # N136: was -12399

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 462.45456  # H7: was 522.8889
$ws.Cells.Item(7, 9).Value = 100.333336  # I7: was 60
$ws.Cells.Item(7, 10).Value = 598.25  # J7: was 655.1429000000001
$ws.Cells.Item(7, 11).Value = 301.000008  # K7: was 180
$ws.Cells.Item(7, 12).Value = 1794.75  # L7: was 1965.4287
$ws.Cells.Item(7, 13).Value = -189.000008  # M7: was -68
$ws.Cells.Item(7, 14).Value = -2018.75  # N7: was -2189.4287
$ws.Cells.Item(80, 8).Value = 2614  # H80: was 2706
$ws.Cells.Item(80, 9).Value = 2921  # I80: was 3666.6667
$ws.Cells.Item(80, 10).Value = 2000  # J80: was 1265
$ws.Cells.Item(80, 11).Value = 8763  # K80: was 11000.0001
$ws.Cells.Item(80, 12).Value = 6000  # L80: was 3795
$ws.Cells.Item(80, 13).Value = -7827  # M80: was -10064.0001
$ws.Cells.Item(80, 14).Value = -7872  # N80: was -5667
$ws.Cells.Item(83, 8).Value = 2614  # H83: was 2706
$ws.Cells.Item(83, 9).Value = 2921  # I83: was 3666.6667
$ws.Cells.Item(83, 10).Value = 2000  # J83: was 1265
$ws.Cells.Item(83, 11).Value = 26289  # K83: was 33000.0003
$ws.Cells.Item(83, 12).Value = 18000  # L83: was 11385
$ws.Cells.Item(83, 13).Value = -21609  # M83: was -28320.0003
$ws.Cells.Item(83, 14).Value = -27360  # N83: was -20745
$ws.Cells.Item(92, 8).Value = 100000410  # H92: was 71428936
$ws.Cells.Item(92, 9).Value = 166666850  # I92: was 125000190
$ws.Cells.Item(92, 10).Value = 750  # J92: was 590
$ws.Cells.Item(92, 11).Value = 500000550  # K92: was 375000570
$ws.Cells.Item(92, 12).Value = 2250  # L92: was 1770
$ws.Cells.Item(92, 13).Value = -499999302  # M92: was -374999322
$ws.Cells.Item(92, 14).Value = -4746  # N92: was -4266
$ws.Cells.Item(113, 8).Value = 15152091  # H113: was 14706453
$ws.Cells.Item(113, 9).Value = 22727826  # I113: was 20833878
$ws.Cells.Item(113, 10).Value = 620.4545000000001  # J113: was 632.5
$ws.Cells.Item(113, 11).Value = 68183478  # K113: was 62501634
$ws.Cells.Item(113, 12).Value = 1861.3635  # L113: was 1897.5
$ws.Cells.Item(113, 13).Value = -68181308  # M113: was -62499464
$ws.Cells.Item(113, 14).Value = -6201.3635  # N113: was -6237.5
$ws.Cells.Item(120, 8).Value = 2000  # H120: was 6557.8
$ws.Cells.Item(120, 10).Value = 0  # J120: was 7697.25
$ws.Cells.Item(120, 12).Value = 0  # L120: was 23091.75
$ws.Cells.Item(120, 14).Value = $null  # N120: was -32767.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(59, 8).Value = 22500  # H59: was 25000
$ws.Cells.Item(59, 10).Value = 22500  # J59: was 25000
$ws.Cells.Item(59, 12).Value = 22500  # L59: was 25000
$ws.Cells.Item(59, 14).Value = -23666  # N59: was -26166
$ws.Cells.Item(69, 8).Value = 33000  # H69: was 0
$ws.Cells.Item(69, 10).Value = 33000  # J69: was 0
$ws.Cells.Item(69, 12).Value = 33000  # L69: was 0
$ws.Cells.Item(69, 14).Value = -34498  # N69: was None
$ws.Cells.Item(72, 8).Value = 33000  # H72: was 0
$ws.Cells.Item(72, 10).Value = 33000  # J72: was 0
$ws.Cells.Item(72, 12).Value = 99000  # L72: was 0
$ws.Cells.Item(72, 14).Value = -106488  # N72: was None
$ws.Cells.Item(113, 8).Value = 4645.724  # H113: was 4793.8213
$ws.Cells.Item(113, 9).Value = 6411.3887  # I113: was 7136.4375
$ws.Cells.Item(113, 10).Value = 1756.4546  # J113: was 1670.3334
$ws.Cells.Item(113, 11).Value = 6411.3887  # K113: was 7136.4375
$ws.Cells.Item(113, 12).Value = 1756.4546  # L113: was 1670.3334
$ws.Cells.Item(113, 13).Value = -4241.3887  # M113: was -4966.4375
$ws.Cells.Item(113, 14).Value = -6096.4546  # N113: was -6010.3334
$ws.Cells.Item(132, 8).Value = 4195.026  # H132: was 4654.8184
$ws.Cells.Item(132, 9).Value = 4641.4  # I132: was 5680.7334
$ws.Cells.Item(132, 10).Value = 3725.158  # J132: was 3799.889
$ws.Cells.Item(132, 11).Value = 13924.2  # K132: was 17042.2002
$ws.Cells.Item(132, 12).Value = 11175.474  # L132: was 11399.667
$ws.Cells.Item(132, 13).Value = -11394.2  # M132: was -14512.2002
$ws.Cells.Item(132, 14).Value = -16235.474  # N132: was -16459.667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2437.5  # H122: was 1770.4286
$ws.Cells.Item(122, 9).Value = 2300  # I122: was 1718.9333
$ws.Cells.Item(122, 10).Value = 2850  # J122: was 1899.1666
$ws.Cells.Item(122, 11).Value = 6900  # K122: was 5156.7999
$ws.Cells.Item(122, 12).Value = 8550  # L122: was 5697.4998
$ws.Cells.Item(122, 13).Value = -4450  # M122: was -2706.7999
$ws.Cells.Item(122, 14).Value = -13450  # N122: was -10597.4998
